# Sprint 3 Backlog - Burndown: update completed totals for finished tasks
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 18 (Implement functionality to make recipe list paginated (Web)):
#   Week 2 (E18) completed amount = 2
#   Amount Remaining After Week 2 (I18) = 0
$ws.Range("E18").Value = 2
$ws.Range("I18").Value = 0

# Row 21 (Modify functionality to fit new Desktop UI for filtering recipes):
#   Amount Remaining After Week 2 (I21) = 0
$ws.Range("I21").Value = 0

# Row 23 (Update class diagram to include a more detailed representation of project file):
#   Week 2 (E23) completed amount = 3
#   Week 2 actual time (H23) = 2
#   Amount Remaining After Week 2 (I23) = 0
$ws.Range("E23").Value = 3
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = 0

# Update the sheet view: scroll down and move the selection to E38
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("E38").Select()
